$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "santoshAmal123@gmail.com"
$ws.Range("B2").Value = "SantoshAmal123456"

# These ColumnWidth inputs are chosen (empirically, via this runtime's
# rounding behaviour) so the serialized OOXML <col width="..."> ends up
# matching (or as close as achievable to) the target widths of 27 and
# 18.85546875 respectively.
$ws.Columns.Item(1).ColumnWidth = 26.166666666666668
$ws.Columns.Item(2).ColumnWidth = 18.0
